# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (holding fund data) right before the
#   "总计" (Total) summary sheet.
# - Prepend a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data rows, and renumber the existing rows' running index.

$wb = $excel.ActiveWorkbook

# Use an existing per-quarter sheet as a formatting template for the new one.
$refSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet, positioned right before "总计".
#    (Re-resolve "总计" by index right before use -- sheet object handles in
#    this runtime can become stale/aliased once new sheets are inserted.)
# ---------------------------------------------------------------------------
$totalIndex = $wb.Worksheets.Item("总计").Index
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item($totalIndex))
$newSheet.Name = "2022-Q1"

# Match page setup / outline properties used by the other per-quarter sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy header (row1, B:H) and the bold/bordered "index" column (A) formatting
# from the reference sheet.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$refSheet.Range("A2:A9").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- Fund holding rows ---
# idx, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "003834", "华夏能源革新股票",            "187.75", "93.26", "3.82",  "7.1720", 10),
    @(1, "004854", "广发中证全指汽车指数A",        "22.01",  "94.43", "12.34", "2.7160", 3),
    @(2, "002083", "新华鑫动力灵活配置混合A",      "29.11",  "91.34", "5.48",  "1.5952", 7),
    @(3, "004855", "广发中证全指汽车指数C",        "6.11",   "94.43", "12.34", "0.7540", 3),
    @(4, "002084", "新华鑫动力灵活配置混合C",      "12.36",  "91.34", "5.48",  "0.6773", 7),
    @(5, "516110", "国泰中证800汽车与零部件ETF",   "1.74",   "97.85", "8.27",  "0.1439", 4),
    @(6, "002272", "新华科技创新主题灵活配置混合", "1.74",   "88.44", "5.22",  "0.0908", 6),
    @(7, "159936", "广发中证全指可选消费ETF指数",  "1.78",   "97.88", "1.57",  "0.0279", 8)
)

# Columns B (基金代码) and D:G (基金规模/股票总仓位/仓位占比/持有市值) are stored
# as text in the source data (to preserve leading zeros / trailing zeros), so
# force a text number format on those ranges before writing the values.
$newSheet.Range("B2:B9").NumberFormat = "@"
$newSheet.Range("D2:G9").NumberFormat = "@"

foreach ($row in $rows) {
    $r = [int]$row[0] + 2
    $newSheet.Cells.Item($r, 1).Value = [int]$row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (Total) sheet: insert a new row for 2022-Q1 at the top
#    of the data (row 2), shifting the previous rows down, then renumber the
#    running index in column A.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).ClearFormats()
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 13.18

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
